# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and "Correspond Handback
# DateTime" (H2) values for the first data row on both the zh-cn and
# de-de worksheets, as part of regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-18 10:35:40"
$wsZhCn.Range("H2").Value = "2016-03-18 10:35:57"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-18 10:35:44"
$wsDeDe.Range("H2").Value = "2016-03-18 10:36:03"
